# Weekly update: insert the new price observation as row 165, pushing the
# existing rows 165-213 down to 166-214 (dimension grows from R213 to R214).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(165).Insert()

$ws.Cells.Item(165, 1).Value = 3
$ws.Cells.Item(165, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(165, 3).Value = "Coquimbo"
$ws.Cells.Item(165, 4).Value = 44463
$ws.Cells.Item(165, 5).Value = 5
$ws.Cells.Item(165, 6).Value = 100112032
$ws.Cells.Item(165, 7).Value = "Zapallo italiano"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 200
$ws.Cells.Item(165, 11).Value = 11000
$ws.Cells.Item(165, 12).Value = 12000
$ws.Cells.Item(165, 13).Value = 11525
$ws.Cells.Item(165, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(165, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(165, 16).Value = 165
$ws.Cells.Item(165, 17).Value = 70
$ws.Cells.Item(165, 18).Value = "Hortaliza"
